$d = $word.ActiveDocument

# Locate the paragraph "Product Development and Platform Architecture"
# (the heading line under the Siege Analytics / PARTNER entry) so we can
# insert the three new bullet paragraphs right after it.
$targetIndex = -1
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text -eq "Product Development and Platform Architecture`r") {
        $targetIndex = $idx
        break
    }
}

if ($targetIndex -eq -1) {
    Write-Output "ERROR: anchor paragraph not found"
} else {
    $newBullets = @(
        "• Conceived and architected redistricting platform incorporating boundary estimation algorithm used by 2,500+ analysts",
        "• Built multi-tenant data warehouse tracking decades of demographic data, enabling discovery of 500,000+ mischaracterized voters",
        "• Platform democratized redistricting analysis, reducing costs by 75% and enabling 200+ smaller organizations to participate"
    )

    $insertIndex = $targetIndex
    foreach ($bulletText in $newBullets) {
        $anchor = $d.Paragraphs.Item($insertIndex)
        $anchor.Range.InsertParagraphAfter()
        $insertIndex = $insertIndex + 1
        $newParagraph = $d.Paragraphs.Item($insertIndex)
        $newParagraph.Range.Text = $bulletText
    }

    Write-Output "Inserted $($newBullets.Length) bullet paragraphs after paragraph $targetIndex"
}
